$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 16.98875427246094
$ws.Range("D2").Value = 177

$ws.Range("C3").Value = 15.77210426330566
$ws.Range("D3").Value = 176

$ws.Range("C4").Value = 15.27690887451172
$ws.Range("D4").Value = 175

$ws.Range("C5").Value = 15.2897834777832
$ws.Range("D5").Value = 179

$ws.Range("C6").Value = 15.15483856201172
$ws.Range("D6").Value = 181
